$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# The target text box ("CasellaDiTesto 43") lives two group-levels deep:
# Gruppo 2 (top-level shape) -> Gruppo 41 (nested group) -> CasellaDiTesto 43.
# This runtime's Shapes/GroupItems traversal flattens the nested group's
# items into the parent group's GroupItems collection, so the textbox is
# item 8 of the outer group's GroupItems.
$grp = $s.Shapes.Item(5)
$items = $grp.GroupItems
$shp = $items.Item(8)

$tr = $shp.TextFrame.TextRange

$newText = "Includere lo scraper all’interno di AWS, per automatizzare il processo di scraping e integrare Amazon Elastic Transcoder per generare la trascrizione di quei talk che non ne hanno una"

# Setting TextRange.Text directly to the merged string would normally keep
# the three original runs split (the host tries to preserve formatting for
# overlapping text), leaving the old underline run boundary around
# "scraper" intact. Route through an unrelated placeholder string first so
# the subsequent assignment has no overlap with the original runs and
# collapses into a single run (inheriting the first run's formatting, with
# no underline), matching the target diff.
$tr.Text = "placeholder text with no overlap zzz"
$tr.Text = $newText

Write-Host "Final text: $($tr.Text)"
Write-Host "Final run count: $($tr.Runs().Count)"
